$p = $ppt.ActivePresentation
$s5 = $p.Slides.Item(5)
$s7 = $p.Slides.Item(7)
$t5 = $s5.TimeLine
$t7 = $s7.TimeLine
$m5 = $t5.MainSequence
$m7 = $t7.MainSequence
Write-Host "done"
